# edit.ps1 - apply the "Miscellaneous tweaks (#5, #4, #3)" commit to before.docx
$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

# ---------------------------------------------------------------------------
# 1) Remove the "This version:" and "Latest version:" paragraphs entirely.
#    (They are paragraphs 3 and 4 at the start of the document.)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p4 = $d.Paragraphs(4)
if ($p3.Range.Text -like "This version:*" -and $p4.Range.Text -like "Latest version:*") {
    $rng = $d.Range($p3.Range.Start, $p4.Range.End)
    $rng.Delete()
}

Write-Output "step1 done"

# ---------------------------------------------------------------------------
# 2) "This Draft Best Practice is published ..." -> "This Best Practice is
#    published ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Draft Best Practice", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Best Practice", $wdReplaceAll) | Out-Null

Write-Output "step2 done"

# ---------------------------------------------------------------------------
# 3) "Status of this Document" paragraphs: move the "it may be updated..."
#    clause into the previous sentence, drop the "should not be cited..."
#    sentence, and turn "Readers are encouraged..." into its own paragraph
#    preceded by a new "latest version" paragraph with a hyperlink.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$insertPoint = $d.Range($p6.Range.End - 1, $p6.Range.End - 1)
$insertPoint.InsertAfter(" It may be updated, replaced or obsoleted by other documents at any time.")

$p7 = $d.Paragraphs(7)
$body = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$body.Text = "The latest version of this document is available at https://www.imfug.com/TR/timecode-in-imf/."

$p7 = $d.Paragraphs(7)
$insertPoint2 = $d.Range($p7.Range.End - 1, $p7.Range.End - 1)
$insertPoint2.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8body = $d.Range($p8.Range.Start, $p8.Range.End - 1)
$p8body.InsertAfter("Readers are encouraged to consult the following for a list of current issues, to which they are invited to contribute.")

# Turn the URL text in the new paragraph into a real hyperlink.
$p7 = $d.Paragraphs(7)
$findRng = $p7.Range.Duplicate()
$found = $findRng.Find.Execute("https://www.imfug.com/TR/timecode-in-imf/", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$d.Hyperlinks.Add($findRng, "https://www.imfug.com/TR/timecode-in-imf/", $null, $null, $findRng.Text) | Out-Null

Write-Output "step3 done"

# ---------------------------------------------------------------------------
# 4) Copyright line: "is © 2018 Hollywood" -> "is © Hollywood" (drop the year)
# ---------------------------------------------------------------------------
$copyrightChar = [char]0x00A9
$d.Content.Find.Execute("is $copyrightChar 2018 Hollywood", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "is $copyrightChar Hollywood", $wdReplaceAll) | Out-Null

Write-Output "step4 done"

# ---------------------------------------------------------------------------
# 5) Number the four "NOTE: " paragraphs as "NOTE 1: " .. "NOTE 4: ".
#    Each paragraph's own Range is searched (not Content) so the global
#    replace doesn't clobber every occurrence with the same number.
# ---------------------------------------------------------------------------
$noteNumber = 1
for ($pi = 1; $pi -le $d.Paragraphs.Count; $pi++) {
    $para = $d.Paragraphs($pi)
    if ($para.Range.Text.StartsWith("NOTE: ")) {
        $rng = $para.Range.Duplicate()
        $rng.Find.Execute("NOTE: ", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "NOTE ${noteNumber}: ", $wdReplaceAll) | Out-Null
        $noteNumber++
    }
}

Write-Output "step5 done"

# ---------------------------------------------------------------------------
# 6) Fix typo in endnote 2: "wwww.imfug.com" -> "www.imfug.com"
# ---------------------------------------------------------------------------
for ($ei = 1; $ei -le $d.Endnotes.Count; $ei++) {
    $en = $d.Endnotes($ei)
    if ($en.Range.Text -like "*wwww.imfug.com*") {
        $en.Range.Text = " https://www.imfug.com/"
    }
}

Write-Output "step6 done"
